$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns: batch_id / semester_id / branch_id -----------------
$ws.Range("E1").Value = "batch_id"
$ws.Range("F1").Value = "semester_id"
$ws.Range("G1").Value = "branch_id"

# G1 is a brand new cell outside the old used range, so it doesn't inherit the
# bold header style the way E1/F1 do (they already existed as header cells) -
# match it to the rest of row 1 explicitly.
$ws.Range("G1").Font.Bold = $true

# --- New data values for rows 2-4 -------------------------------------------
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 1

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 1

# --- Column widths: split the old shared E:F width into individual ones -----
$ws.Columns.Item(5).ColumnWidth = 9.43877551020408
$ws.Columns.Item(6).ColumnWidth = 11.8061224489796

# --- Match the selection left behind by the edit -----------------------------
[void]$ws.Range("G4").Select()
